# Scheduled runner update: refresh Kraken market-price-derived figures
# (currentAveragePrice*, Leve*Price*, Leve*Profit*) across the per-job
# profit sheets. Values below mirror a re-pull of Kraken's NA Leve profit
# data; unaffected rows/columns are left untouched.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 3412.0908
$ws.Range("J41").Value = 5062.2856
$ws.Range("L41").Value = 5062.2856
$ws.Range("N41").Value = -5942.2856

$ws.Range("H74").Value = 4300
$ws.Range("I74").Value = 5000
$ws.Range("K74").Value = 5000
$ws.Range("M74").Value = -4064

$ws.Range("H77").Value = 4300
$ws.Range("I77").Value = 5000
$ws.Range("K77").Value = 25000
$ws.Range("M77").Value = -20320

$ws.Range("H106").Value = 5000
$ws.Range("J106").Value = 5000
$ws.Range("L106").Value = 5000
$ws.Range("N106").Value = -6262

$ws.Range("H127").Value = 2524.4167
$ws.Range("I127").Value = 679.6667
$ws.Range("J127").Value = 4369.1665
$ws.Range("K127").Value = 2039.0001
$ws.Range("L127").Value = 13107.4995
$ws.Range("M127").Value = 2920.9999
$ws.Range("N127").Value = -23027.4995

$ws.Range("H132").Value = 6779.7334
$ws.Range("J132").Value = 9499.799999999999
$ws.Range("L132").Value = 28499.4
$ws.Range("N132").Value = -33559.39999999999

$ws.Range("H135").Value = 1413.9231
$ws.Range("I135").Value = 894.7778
$ws.Range("K135").Value = 8053.000199999999
$ws.Range("M135").Value = -5518.000199999999

$ws.Range("H137").Value = 1890.8214
$ws.Range("I137").Value = 1673.5883
$ws.Range("K137").Value = 5020.7649
$ws.Range("M137").Value = -2470.7649

$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").Value = $null

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6088.4585
$ws.Range("I32").Value = 4987.15
$ws.Range("K32").Value = 4987.15
$ws.Range("M32").Value = -4700.15

$ws.Range("H97").Value = 1553.625
$ws.Range("I97").Value = 571.5
$ws.Range("J97").Value = 4500
$ws.Range("K97").Value = 571.5
$ws.Range("L97").Value = 4500
$ws.Range("M97").Value = -75.5
$ws.Range("N97").Value = -5492

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1790
$ws.Range("I107").Value = 1790
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1790
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 130
$ws.Range("N107").Value = $null

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 635.5
$ws.Range("I14").Value = 1000
$ws.Range("J14").Value = 271
$ws.Range("K14").Value = 1000
$ws.Range("L14").Value = 271
$ws.Range("M14").Value = -830
$ws.Range("N14").Value = -611

$ws.Range("H16").Value = 1100
$ws.Range("I16").Value = 1100
$ws.Range("K16").Value = 1100
$ws.Range("M16").Value = -813

$ws.Range("H107").Value = 387.4
$ws.Range("I107").Value = 216.33333
$ws.Range("K107").Value = 216.33333
$ws.Range("M107").Value = 1703.66667

$ws.Range("H113").Value = 1100
$ws.Range("I113").Value = 1100
$ws.Range("K113").Value = 1100
$ws.Range("M113").Value = 1070

$ws.Range("H132").Value = 7655.75
$ws.Range("I132").Value = 7640.857
$ws.Range("K132").Value = 22922.571
$ws.Range("M132").Value = -20392.571

$ws.Range("H134").Value = 1075
$ws.Range("I134").Value = 1075
$ws.Range("K134").Value = 3225
$ws.Range("M134").Value = -690

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 165
$ws.Range("I60").Value = 165
$ws.Range("K60").Value = 495
$ws.Range("M60").Value = -244

$ws.Range("H129").Value = 2005.1666
$ws.Range("I129").Value = 1009.1429
$ws.Range("J129").Value = 3399.6
$ws.Range("K129").Value = 3027.4287
$ws.Range("L129").Value = 10198.8
$ws.Range("M129").Value = 1972.5713
$ws.Range("N129").Value = -20198.8

$ws.Range("H131").Value = 2689.9167
$ws.Range("I131").Value = 2715.6365
$ws.Range("J131").Value = 2668.1538
$ws.Range("K131").Value = 8146.9095
$ws.Range("L131").Value = 8004.4614
$ws.Range("M131").Value = -3106.9095
$ws.Range("N131").Value = -18084.4614

$ws.Range("H140").Value = 1341
$ws.Range("I140").Value = 1341
$ws.Range("K140").Value = 4023
$ws.Range("M140").Value = 1157

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 6000000
$ws.Range("I12").Value = 6000000
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 6000000
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -5999860
$ws.Range("N12").Value = $null

$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = $null
$ws.Range("N24").Value = $null

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 236.66667
$ws.Range("I16").Value = 236.66667
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 236.66667
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -66.66667000000001
$ws.Range("N16").Value = $null

$ws.Range("H55").Value = 1942.7142
$ws.Range("I55").Value = 899.75
$ws.Range("J55").Value = 3333.3333
$ws.Range("K55").Value = 899.75
$ws.Range("L55").Value = 3333.3333
$ws.Range("M55").Value = -726.75
$ws.Range("N55").Value = -3679.3333

$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").Value = $null

$ws.Range("H132").Value = 28931.25
$ws.Range("I132").Value = 29908.334
$ws.Range("J132").Value = 26000
$ws.Range("K132").Value = 89725.00199999999
$ws.Range("L132").Value = 78000
$ws.Range("M132").Value = -87195.00199999999
$ws.Range("N132").Value = -83060

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 30003
$ws.Range("I18").Value = 30003
$ws.Range("K18").Value = 30003
$ws.Range("M18").Value = -29830

$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = $null
$ws.Range("N20").Value = $null

$ws.Range("H21").Value = 1894757.2
$ws.Range("J21").Value = 25015
$ws.Range("L21").Value = 25015
$ws.Range("N21").Value = -25485

$ws.Range("H29").Value = 11401
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 11401
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 11401
$ws.Range("M29").Value = $null
$ws.Range("N29").Value = -11981

$ws.Range("H31").Value = 26016.5
$ws.Range("I31").Value = 26016.5
$ws.Range("K31").Value = 26016.5
$ws.Range("M31").Value = -25668.5

$ws.Range("H35").Value = 1894757.2
$ws.Range("J35").Value = 25015
$ws.Range("L35").Value = 25015
$ws.Range("N35").Value = -25595

$ws.Range("H113").Value = 463.9375
$ws.Range("I113").Value = 430.2143
$ws.Range("K113").Value = 1290.6429
$ws.Range("M113").Value = 879.3571000000002

$ws.Range("H128").Value = 75000
$ws.Range("J128").Value = 75000
$ws.Range("L128").Value = 75000
$ws.Range("N128").Value = -84960

$ws.Range("H132").Value = 8925.666999999999
$ws.Range("I132").Value = 6389
$ws.Range("K132").Value = 19167
$ws.Range("M132").Value = -16637
